$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update case with 380 kV: replace computed values for rows 2-25
# across columns B, D, E, F, G, I, M, N (pl_mw.xlsx line-results sheet)
$ws.Range("B2").Value = 0.123079927955402
$ws.Range("D2").Value = 0.314049879459219
$ws.Range("E2").Value = 0.1768047495715734
$ws.Range("F2").Value = 4.749763998417563
$ws.Range("G2").Value = 0.002376493410274945
$ws.Range("I2").Value = 0.3320621333132454
$ws.Range("M2").Value = 2.292532481325765
$ws.Range("N2").Value = 1.043668644564235
$ws.Range("B3").Value = 0.1092979331872783
$ws.Range("D3").Value = 0.2807962535968045
$ws.Range("E3").Value = 0.1543910560365447
$ws.Range("F3").Value = 4.302816729529354
$ws.Range("G3").Value = 0.002393512444321289
$ws.Range("I3").Value = 0.3389284839788793
$ws.Range("M3").Value = 2.021132161612769
$ws.Range("N3").Value = 1.046502041668859
$ws.Range("B4").Value = 0.1008218889714954
$ws.Range("D4").Value = 0.2607624217283728
$ws.Range("E4").Value = 0.1406496207408239
$ws.Range("F4").Value = 4.03393331036699
$ws.Range("G4").Value = 0.00240442285053875
$ws.Range("I4").Value = 0.3433865055670982
$ws.Range("M4").Value = 1.856072850322704
$ws.Range("N4").Value = 1.048757236476902
$ws.Range("B5").Value = 0.09736471753583942
$ws.Range("D5").Value = 0.2526883274542797
$ws.Range("E5").Value = 0.1350528360446432
$ws.Range("F5").Value = 3.925670506065671
$ws.Range("G5").Value = 0.002408985852717006
$ws.Range("I5").Value = 0.3452640908047311
$ws.Range("M5").Value = 1.789173779688667
$ws.Range("N5").Value = 1.049806739674935
$ws.Range("B6").Value = 0.09679047958253761
$ws.Range("D6").Value = 0.2513528659874567
$ws.Range("E6").Value = 0.1341236096999907
$ws.Range("F6").Value = 3.907770219238387
$ws.Range("G6").Value = 0.002409750628937506
$ws.Range("I6").Value = 0.3455795420785019
$ws.Range("M6").Value = 1.778086260533911
$ws.Range("N6").Value = 1.049988916084473
$ws.Range("B7").Value = 0.1007752764710261
$ws.Range("D7").Value = 0.2606531765314912
$ws.Range("E7").Value = 0.1405741314336524
$ws.Range("F7").Value = 4.032468048268242
$ws.Range("G7").Value = 0.002404483913899136
$ws.Range("I7").Value = 0.343411580681126
$ws.Range("M7").Value = 1.855169196098529
$ws.Range("N7").Value = 1.048770860916775
$ws.Range("B8").Value = 0.1183309839436646
$ws.Range("D8").Value = 0.3025001936304932
$ws.Range("E8").Value = 0.1690708372857088
$ws.Range("F8").Value = 4.59445377178281
$ws.Range("G8").Value = 0.002382266630113423
$ws.Range("I8").Value = 0.3343794363853387
$ws.Range("M8").Value = 2.198604562861959
$ws.Range("N8").Value = 1.044539131261629
$ws.Range("B9").Value = 0.1526335022417413
$ws.Range("D9").Value = 0.3879284170784842
$ws.Range("E9").Value = 0.2252249949888068
$ws.Range("F9").Value = 5.744530977641318
$ws.Range("G9").Value = 0.00234230273841632
$ws.Range("I9").Value = 0.3185860856400566
$ws.Range("M9").Value = 2.886251940269631
$ws.Range("N9").Value = 1.04029350483421
$ws.Range("B10").Value = 0.1777432096085079
$ws.Range("D10").Value = 0.4532134133729357
$ws.Range("E10").Value = 0.2668092318746034
$ws.Range("F10").Value = 6.624702995414623
$ws.Range("G10").Value = 0.002315065032966163
$ws.Range("I10").Value = 0.3081494023467606
$ws.Range("M10").Value = 3.402510883514594
$ws.Range("N10").Value = 1.039596101849781
$ws.Range("B11").Value = 0.1891428895084459
$ws.Range("D11").Value = 0.4835700066855679
$ws.Range("E11").Value = 0.2858363802429835
$ws.Range("F11").Value = 7.034146820574165
$ws.Range("G11").Value = 0.002303118397442461
$ws.Range("I11").Value = 0.3036541772144457
$ws.Range("M11").Value = 3.640332165465765
$ws.Range("N11").Value = 1.039794476428213
$ws.Range("B12").Value = 0.1934560386771409
$ws.Range("D12").Value = 0.4951694077117281
$ws.Range("E12").Value = 0.2930606355047445
$ws.Range("F12").Value = 7.190614112860715
$ws.Range("G12").Value = 0.002298656954922684
$ws.Range("I12").Value = 0.3019882179170086
$ws.Range("M12").Value = 3.730865017624069
$ws.Range("N12").Value = 1.039942876553752
$ws.Range("B13").Value = 0.1925272936593103
$ws.Range("D13").Value = 0.4926664877883979
$ws.Range("E13").Value = 0.291503860275057
$ws.Range("F13").Value = 7.156850975244083
$ws.Range("G13").Value = 0.002299615049206621
$ws.Range("I13").Value = 0.3023453986127027
$ws.Range("M13").Value = 3.711345229054615
$ws.Range("N13").Value = 1.039907670603455
$ws.Range("B14").Value = 0.1894978104235463
$ws.Range("D14").Value = 0.4845221560613027
$ws.Range("E14").Value = 0.2864303222832376
$ws.Range("F14").Value = 7.046990320650991
$ws.Range("G14").Value = 0.002302750107561651
$ws.Range("I14").Value = 0.3035163906564442
$ws.Range("M14").Value = 3.647770553082296
$ws.Range("N14").Value = 1.039805221724905
$ws.Range("B15").Value = 0.1876416764396822
$ws.Range("D15").Value = 0.4795473515040101
$ws.Range("E15").Value = 0.2833252189298179
$ws.Range("F15").Value = 6.979886031882074
$ws.Range("G15").Value = 0.002304678517939509
$ws.Range("I15").Value = 0.3042383818420014
$ws.Range("M15").Value = 3.608892562042968
$ws.Range("N15").Value = 1.039751986303813
$ws.Range("B16").Value = 0.1769977216430618
$ws.Range("D16").Value = 0.4512435715029142
$ws.Range("E16").Value = 0.2655682533169994
$ws.Range("F16").Value = 6.59813695115497
$ws.Range("G16").Value = 0.002315854592897788
$ws.Range("I16").Value = 0.308448254558785
$ws.Range("M16").Value = 3.387032527645118
$ws.Range("N16").Value = 1.039593439288623
$ws.Range("B17").Value = 0.1704618784255132
$ws.Range("D17").Value = 0.434055188383752
$ws.Range("E17").Value = 0.25470540650619
$ws.Range("F17").Value = 6.366346205026503
$ws.Range("G17").Value = 0.002322823512681115
$ws.Range("I17").Value = 0.3110955253145953
$ws.Range("M17").Value = 3.251722943464245
$ws.Range("N17").Value = 1.039627623963057
$ws.Range("B18").Value = 0.1667005055892474
$ws.Range("D18").Value = 0.4242302381405239
$ws.Range("E18").Value = 0.248467558086702
$ws.Range("F18").Value = 6.233870562101345
$ws.Range("G18").Value = 0.002326873719914264
$ws.Range("I18").Value = 0.3126419307778239
$ws.Range("M18").Value = 3.174173220501302
$ws.Range("N18").Value = 1.039695871109501
$ws.Range("B19").Value = 0.1654266152931285
$ws.Range("D19").Value = 0.4209139390295036
$ws.Range("E19").Value = 0.2463571820807857
$ws.Range("F19").Value = 6.189158085657766
$ws.Range("G19").Value = 0.002328252281894693
$ws.Range("I19").Value = 0.3131695994379697
$ws.Range("M19").Value = 3.147962345755673
$ws.Range("N19").Value = 1.03972734932826
$ws.Range("B20").Value = 0.1711578531401301
$ws.Range("D20").Value = 0.4358785050330312
$ws.Range("E20").Value = 0.2558606994092116
$ws.Range("F20").Value = 6.39093251027009
$ws.Range("G20").Value = 0.002322077335447833
$ws.Range("I20").Value = 0.3108112593494798
$ws.Range("M20").Value = 3.266097898637184
$ws.Range("N20").Value = 1.039618962853979
$ws.Range("B21").Value = 0.1903877454561496
$ws.Range("D21").Value = 0.486911443256588
$ws.Range("E21").Value = 0.2879199985655703
$ws.Range("F21").Value = 7.079219546360264
$ws.Range("G21").Value = 0.002301827579918527
$ws.Range("I21").Value = 0.3031714573127502
$ws.Range("M21").Value = 3.666430688752484
$ws.Range("N21").Value = 1.039833331560644
$ws.Range("B22").Value = 0.2029340735935818
$ws.Range("D22").Value = 0.5208754614550912
$ws.Range("E22").Value = 0.3089855113396567
$ws.Range("F22").Value = 7.537390591756491
$ws.Range("G22").Value = 0.002288956761662408
$ws.Range("I22").Value = 0.2983898912704426
$ws.Range("M22").Value = 3.930864699697736
$ws.Range("N22").Value = 1.040400161853981
$ws.Range("B23").Value = 0.1962399543991182
$ws.Range("D23").Value = 0.5026890735056782
$ws.Range("E23").Value = 0.2977309966580464
$ws.Range("F23").Value = 7.292052399798479
$ws.Range("G23").Value = 0.002295793353588052
$ws.Range("I23").Value = 0.3009225589881588
$ws.Range("M23").Value = 3.789459204654207
$ws.Range("N23").Value = 1.040058875075019
$ws.Range("B24").Value = 0.1708432148711267
$ws.Range("D24").Value = 0.4350540071096702
$ws.Range("E24").Value = 0.2553383688454716
$ws.Range("F24").Value = 6.379814608741867
$ws.Range("G24").Value = 0.002322414546157039
$ws.Range("I24").Value = 0.3109396998212492
$ws.Range("M24").Value = 3.259598232569914
$ws.Range("N24").Value = 1.039622727116566
$ws.Range("B25").Value = 0.1433687920817732
$ws.Range("D25").Value = 0.3644104759966922
$ws.Range("E25").Value = 0.2099885130927532
$ws.Range("F25").Value = 5.427668866972226
$ws.Range("G25").Value = 0.002352735617854966
$ws.Range("I25").Value = 0.3226535088687772
$ws.Range("M25").Value = 2.698467210663409
$ws.Range("N25").Value = 1.041013246675064
